$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing all existing data (the ticker list)
# down by one row.
$ws.Rows.Item(1).Insert()

# Set the new header cell.
$ws.Range("A1").Value = "Tickers"

# The row that used to be row 1 (now row 2) ends up with an explicit
# 14pt custom row height after the insert, matching the authored file.
$ws.Rows.Item(2).RowHeight = 14

# Update the selection to match the authored workbook (B2 selected).
$ws.Range("B2").Select()
